# Apply the commit "define interactions parameters in excel":
#  - add a new "parameters" worksheet (after "units") listing the
#    interaction-distance parameters and their numeric values
#  - make the new sheet the active/selected sheet

$wb = $excel.ActiveWorkbook

$units = $wb.Worksheets.Item("units")

# Add the new worksheet right after the last existing sheet ("units")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "parameters"

# Header row
$ws.Range("A1").Value = "parameter"
$ws.Range("B1").Value = "value"

# Reuse the existing bold+italic header style already used on the other
# sheets (e.g. units!A1) instead of re-toggling Font.Bold/Italic, which
# would otherwise create a brand-new (unused) style entry.
$units.Range("A1").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Parameter rows
$data = @(
    @("melee_distance", 3.5),
    @("melee_height_difference_threshold", 2),
    @("archer_distance", 4.5),
    @("archer_distance_height_gain", 0.5),
    @("siege_distance", 11),
    @("siege_distance_height_gain", 0.5),
    @("flier_distance", 10),
    @("flier_distance_height_gain", 0)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# Autofit the columns to their content
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Match the selection left behind in the source file and make this the
# active (tab-selected) sheet
$ws.Range("G10").Select() | Out-Null
$ws.Activate()
